# RB-Blessing / Paarweiser Vergleich — "Formatierung Anforderungsliste und Paarweiser Vergleich"
#
# Replaces the requirement "Getriebemotor" (row 8 / lfd. Nr. 5) with
# "Nothalt mit Fußpedal" and updates that row's pairwise-comparison scores,
# then applies the accompanying view/formatting tweaks (zoom, column widths,
# row heights, a colour-scale on the normalised-factor column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content change: requirement text for lfd. Nr. 5 -----------------------
$ws.Range("C8").Value = "Nothalt mit Fußpedal"

# --- Pairwise-comparison scores for that requirement (row 8) ---------------
# Symmetric cells elsewhere in the sheet (formulas like "=2-H4") recompute
# automatically once these change.
$ws.Range("I8").Value = 2
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 2
$ws.Range("S8").Value = 2
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 2

# --- Row heights -------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 267
$ws.Rows.Item(17).RowHeight = 30.75
$ws.Rows.Item(23).RowHeight = 30.75

# --- Column widths (Q and W get a touch wider) ------------------------
$ws.Columns.Item(17).ColumnWidth = 4.643
$ws.Columns.Item(23).ColumnWidth = 4.072

# --- View: zoom to 85%, scroll/select so AC11 is the active cell ----------
$excel.ActiveWindow.Zoom = 85
$ws.Range("AC11").Select() | Out-Null

# --- Conditional formatting: 3-colour scale on the normalised-factor column
$ws.Range("Y4:Y23").FormatConditions.AddColorScale(3) | Out-Null
